$d = $word.ActiveDocument

# The last two paragraphs of the document are empty paragraphs (with
# ind left=720 hanging=720) left over after the final reference's
# hyperlink. Remove them, leaving the hyperlink paragraph as the last
# paragraph of the body before the section break.

$count = $d.Paragraphs.Count
$d.Paragraphs.Item($count).Range.Delete()
$d.Paragraphs.Item($count - 1).Range.Delete()
